$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 7 (SignUp data row) content updates.
# A7 keeps the same shared-string slot but the text itself changed.
$ws.Range("A7").Value = "venkat"
# B7 now reuses the existing "prasad" string instead of "Shetty".
$ws.Range("B7").Value = "prasad"
# C7 gets a new email address.
$ws.Range("C7").Value = "prasad8985117@yopmail.com"
# D7/E7 both reuse the same "Ganesh979612" text (same shared string slot).
$ws.Range("D7").Value = "Ganesh979612"
$ws.Range("E7").Value = "Ganesh979612"

# Move the selection to F7, matching the saved view state.
$ws.Activate()
$ws.Range("F7").Select()
